$rows = @(
    @(8.906151666666668,26.718455,0.1245005002255258,0.1245005002255258,3.0,1.0,1.635346666666667,4.90604,0.02683720313876748,0.02683720313876747,14.56464544091111,131.0818089682,0.003341245215430601,0.003341245215430601),
    @(8.906151666666668,26.718455,0.1245005002255258,0.1245005002255258,3.0,1.0,45.91636366666668,137.749091,0.7535202194331003,0.7535202194331002,408.9380987971563,3680.442889174406,0.09381364424946895,0.09381364424946892),
    @(8.906151666666668,26.718455,0.1245005002255258,0.1245005002255258,3.0,1.0,0.8399643333333332,2.519893,0.01378441275019327,0.01378441275019327,7.480849747257222,67.327647725315,0.001716166282714177,0.001716166282714177),
    @(8.906151666666668,26.718455,0.1245005002255258,0.1245005002255258,3.0,1.0,7.591029666666667,22.773089,0.1245742015128762,0.1245742015128762,67.60686151749945,608.461753657495,0.01550955040354854,0.01550955040354853),
    @(8.906151666666668,26.718455,0.1245005002255258,0.1245005002255258,3.0,1.0,4.780457,14.341371,0.07845070297336118,0.07845070297336117,42.57547507797834,383.1792757018051,0.00976715176322761,0.009767151763227606),
    @(8.906151666666668,26.718455,0.1245005002255258,0.1245005002255258,2.0,0.6666666666666666,0.172647,0.517941,0.002833260191701732,0.002833260191701732,1.537620366795,13.838583301155,0.0003527423111359347,0.0003527423111359346),
    @(0.4515893333333333,1.354768,0.006312838586270617,0.006312838586270617,3.0,1.0,1.635346666666667,4.90604,0.02683720313876748,0.02683720313876747,0.7385051109688889,6.64654599872,0.0001694189315219942,0.0001694189315219942),
    @(0.4515893333333333,1.354768,0.006312838586270617,0.006312838586270617,3.0,1.0,45.91636366666668,137.749091,0.7535202194331003,0.7535202194331002,20.73534005732089,186.618060515888,0.004756851516772378,0.004756851516772377),
    @(0.4515893333333333,1.354768,0.006312838586270617,0.006312838586270617,3.0,1.0,0.8399643333333332,2.519893,0.01378441275019327,0.01378441275019327,0.3793189333137778,3.413870399823999,0.00008701877269850074,0.00008701877269850072),
    @(0.4515893333333333,1.354768,0.006312838586270617,0.006312838586270617,3.0,1.0,7.591029666666667,22.773089,0.1245742015128762,0.1245742015128762,3.428028026483556,30.852252238352,0.0007864168261643362,0.0007864168261643361),
    @(0.4515893333333333,1.354768,0.006312838586270617,0.006312838586270617,3.0,1.0,4.780457,14.341371,0.07845070297336118,0.07845070297336117,2.158803389658667,19.429230506928,0.0004952466248502894,0.0004952466248502893),
    @(0.4515893333333333,1.354768,0.006312838586270617,0.006312838586270617,2.0,0.6666666666666666,0.172647,0.517941,0.002833260191701732,0.002833260191701732,0.077965543632,0.701689892688,0.00001788591426311918,0.00001788591426311917),
    @(7.781650666666667,23.344952,0.1087809232135948,0.1087809232135948,3.0,1.0,1.635346666666667,4.90604,0.02683720313876748,0.02683720313876747,12.72569647889778,114.53126831008,0.00291937573390591,0.00291937573390591),
    @(7.781650666666667,23.344952,0.1087809232135948,0.1087809232135948,3.0,1.0,45.91636366666668,137.749091,0.7535202194331003,0.7535202194331002,357.3051019376259,3215.745917438632,0.08196862513004319,0.08196862513004317),
    @(7.781650666666667,23.344952,0.1087809232135948,0.1087809232135948,3.0,1.0,0.8399643333333332,2.519893,0.01378441275019327,0.01378441275019327,6.536309014459555,58.82678113013599,0.001499481144923271,0.001499481144923271),
    @(7.781650666666667,23.344952,0.1087809232135948,0.1087809232135948,3.0,1.0,7.591029666666667,22.773089,0.1245742015128762,0.1245742015128762,59.07074106630311,531.636669596728,0.01355129664916707,0.01355129664916707),
    @(7.781650666666667,23.344952,0.1087809232135948,0.1087809232135948,3.0,1.0,4.780457,14.341371,0.07845070297336118,0.07845070297336117,37.19984640102133,334.798617609192,0.008533939896197735,0.008533939896197735),
    @(7.781650666666667,23.344952,0.1087809232135948,0.1087809232135948,2.0,0.6666666666666666,0.172647,0.517941,0.002833260191701732,0.002833260191701732,1.343478642648,12.091307783832,0.000308204659357641,0.0003082046593576409),
    @(3.892567333333333,11.677702,0.05441481329981927,0.05441481329981927,3.0,1.0,1.635346666666667,4.90604,0.02683720313876748,0.02683720313876747,6.365697013342222,57.29127312008,0.001460341398285356,0.001460341398285356),
    @(3.892567333333333,11.677702,0.05441481329981927,0.05441481329981927,3.0,1.0,45.91636366666668,137.749091,0.7535202194331003,0.7535202194331002,178.7325372743203,1608.592835468882,0.04100266205809101,0.04100266205809099),
    @(3.892567333333333,11.677702,0.05441481329981927,0.05441481329981927,3.0,1.0,0.8399643333333332,2.519893,0.01378441275019327,0.01378441275019327,3.269617725098444,29.426559525886,0.000750076246249415,0.0007500762462494148),
    @(3.892567333333333,11.677702,0.05441481329981927,0.05441481329981927,3.0,1.0,7.591029666666667,22.773089,0.1245742015128762,0.1245742015128762,29.54859410683089,265.937346961478,0.006778681917297221,0.006778681917297219),
    @(3.892567333333333,11.677702,0.05441481329981927,0.05441481329981927,3.0,1.0,4.780457,14.341371,0.07845070297336118,0.07845070297336117,18.60825075660467,167.474256809442,0.004268880355535025,0.004268880355535024),
    @(3.892567333333333,11.677702,0.05441481329981927,0.05441481329981927,2.0,0.6666666666666666,0.172647,0.517941,0.002833260191701732,0.002833260191701732,0.672040072398,6.048360651582,0.0001541713243612599,0.0001541713243612599),
    @(43.49559133333333,130.486774,0.6080317382054886,0.6080317382054886,3.0,1.0,1.635346666666667,4.90604,0.02683720313876748,0.02683720313876747,71.13037030166223,640.17333271496,0.01631787127303858,0.01631787127303858),
    @(43.49559133333333,130.486774,0.6080317382054886,0.6080317382054886,3.0,1.0,45.91636366666668,137.749091,0.7535202194331003,0.7535202194331002,1997.159389558049,17974.43450602244,0.4581642087948892,0.4581642087948891),
    @(43.49559133333333,130.486774,0.6080317382054886,0.6080317382054886,3.0,1.0,0.8399643333333332,2.519893,0.01378441275019327,0.01378441275019327,36.53474537724244,328.812708395182,0.008381360444641913,0.008381360444641911),
    @(43.49559133333333,130.486774,0.6080317382054886,0.6080317382054886,3.0,1.0,7.591029666666667,22.773089,0.1245742015128762,0.1245742015128762,330.1763241805428,2971.586917624886,0.07574506828143492,0.0757450682814349),
    @(43.49559133333333,130.486774,0.6080317382054886,0.6080317382054886,3.0,1.0,4.780457,14.341371,0.07845070297336118,0.07845070297336117,207.9288040585727,1871.359236527154,0.04770051729233529,0.04770051729233529),
    @(43.49559133333333,130.486774,0.6080317382054886,0.6080317382054886,2.0,0.6666666666666666,0.172647,0.517941,0.002833260191701732,0.002833260191701732,7.509383356926,67.58445021233399,0.00172271211914882,0.00172271211914882),
    @(7.007517,21.022551,0.09795918646930096,0.09795918646930096,3.0,1.0,1.635346666666667,4.90604,0.02683720313876748,0.02683720313876747,11.45971956756,103.13747610804,0.002628950586585032,0.002628950586585032),
    @(7.007517,21.022551,0.09795918646930096,0.09795918646930096,3.0,1.0,45.91636366666668,137.749091,0.7535202194331003,0.7535202194331002,321.759698972349,2895.837290751142,0.07381422768383565,0.07381422768383564),
    @(7.007517,21.022551,0.09795918646930096,0.09795918646930096,3.0,1.0,0.8399643333333332,2.519893,0.01378441275019327,0.01378441275019327,5.886064345226999,52.97457910704299,0.001350309858965992,0.001350309858965992),
    @(7.007517,21.022551,0.09795918646930096,0.09795918646930096,3.0,1.0,7.591029666666667,22.773089,0.1245742015128762,0.1245742015128762,53.194269436671,478.748424930039,0.01220318743526411,0.01220318743526411),
    @(7.007517,21.022551,0.09795918646930096,0.09795918646930096,3.0,1.0,4.780457,14.341371,0.07845070297336118,0.07845070297336117,33.499133695269,301.492203257421,0.007684967041215232,0.00768496704121523),
    @(7.007517,21.022551,0.09795918646930096,0.09795918646930096,2.0,0.6666666666666666,0.172647,0.517941,0.002833260191701732,0.002833260191701732,1.209826787499,10.888441087491,0.0002775438634349574,0.0002775438634349573)
)

$ws = $excel.ActiveWorkbook.ActiveSheet
$r = 2
foreach ($row in $rows) {
    $c = 7
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
